# Generate Report for Handoff
#
# The handoff-status report tool reran: the "Ready for handoff" status
# (zh-cn / de-de, on both the per-language sheets and the Overview roll-up)
# advanced to "In Translation", and the associated "Latest HO Xliff
# Generate Date" / "Latest Handoff Datetime" timestamps were refreshed to
# the new run's timestamps. The Status / Handoff-Datetime columns also
# shrank to fit the new (shorter) values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet (row 2): zh-cn / de-de status columns (E, F) + Latest HO Xliff Generate Date (G) ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-13 06:56:48"

# --- zh-cn sheet (row 2): Status column (C) + Latest Handoff Datetime (H) ---
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2016-10-13 06:56:38"

# --- de-de sheet (row 2): Status column (C) + Latest Handoff Datetime (H) ---
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2016-10-13 06:56:48"

# --- Column width adjustments: the Status/Handoff-Datetime columns narrow
#     from ~17.22 chars down to ~13.41 chars now that "In Translation" is
#     shorter than "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
